# Updates the "cryptos" price/volume table (rows 2-51) with refreshed
# market data, and rotates the EnergySwap/Decentraland/PancakeSwap rows
# (44-46) so each keeps its link/price/volume but shifts to a new name,
# matching the upstream GitHub Actions scrape commit.
#
# Price cells in column D frequently look like plain numbers (e.g. "1.007",
# "0.6030", "13.14"). Assigning such a string straight to .Value lets Excel
# auto-convert it to a floating point number (losing the original text
# formatting / exact digits, e.g. "13.14" -> 13.1400000000001). To keep
# these as text exactly as scraped, such cells are first marked as Text
# ("@") format, given their value, then restored to the default "Normal"
# style so no stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.274.01'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '1.807.72'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.007'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5214'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3827'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07941'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.100'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.338'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.006'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = '1.816.49'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.349'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001091'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06610'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.006'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.955'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.79%  '
$ws.Range('D23').Value = '28.315.40'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.238'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('D28').Value = '2.020.92'
$ws.Range('E28').Value = '  -0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.372'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.70'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.35%  '
$ws.Range('E31').Value = '  +1.22%  '
$ws.Range('E32').Value = '  -4.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.676'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.604'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07169'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.10'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2172'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02316'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.703'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.89%  '
$ws.Range('E40').Value = '  -2.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6188'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.170'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.386'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6030'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.778'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.80%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.14'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '125.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.202'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.934'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06847'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.92'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.24%  '
